$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 23.63000000000025
$ws.Range("H2").Value = 0.0005937127017192623
$ws.Range("I2").Value = 0.0005937127017192623
$ws.Range("L2").Value = 40.00446661820627
$ws.Range("M2").Value = '[15.483826838742843, 64.52510639766969]'
$ws.Range("N2").Value = 0.001974741821843429
$ws.Range("O2").Value = 0.001974741821843429
$ws.Range("P2").Value = 1.150973885098963
$ws.Range("Q2").Value = '[0.47171060864711656, 1.8302371615508104]'
$ws.Range("R2").Value = 0.001369965683592156
$ws.Range("S2").Value = 0.001369965683592156
$ws.Range("T2").Value = 63.63634727579967
$ws.Range("U2").Value = '[50.13207470340309, 77.14061984819625]'
$ws.Range("V2").Value = [double]"2.602140725116442e-12"
$ws.Range("W2").Value = [double]"2.602140725116442e-12"
$ws.Range("X2").Value = 19.30138138138159
$ws.Range("Y2").Value = 16.74678678678697
$ws.Range("Z2").Value = 21.85597597597621

$ws.Range("F3").Value = 23.63000000000025
$ws.Range("H3").Value = [double]"1.949625142982914e-05"
$ws.Range("I3").Value = [double]"1.949625142982914e-05"
$ws.Range("L3").Value = 50.08368576029952
$ws.Range("M3").Value = '[25.765203416512122, 74.40216810408693]'
$ws.Range("N3").Value = 0.000146709781816412
$ws.Range("O3").Value = 0.000146709781816412
$ws.Range("P3").Value = 1.163552834662886
$ws.Range("Q3").Value = '[0.6226580034141929, 1.7044476659115784]'
$ws.Range("R3").Value = [double]"8.159274839791841e-05"
$ws.Range("S3").Value = [double]"8.159274839791841e-05"
$ws.Range("T3").Value = 56.78623985073335
$ws.Range("U3").Value = '[43.347214872581134, 70.22526482888556]'
$ws.Range("V3").Value = [double]"6.303668698137699e-11"
$ws.Range("W3").Value = [double]"6.303668698137699e-11"
$ws.Range("X3").Value = 19.25407407407429
$ws.Range("Y3").Value = 17.21985985986005
$ws.Range("Z3").Value = 21.28828828828852

$ws.Range("F4").Value = 23.63000000000025
$ws.Range("H4").Value = 0.0019095733778145
$ws.Range("I4").Value = 0.0019095733778145
$ws.Range("L4").Value = 37.24086796399508
$ws.Range("M4").Value = '[11.350469615988985, 63.13126631200117]'
$ws.Range("N4").Value = 0.005797800840361855
$ws.Range("O4").Value = 0.005797800840361855
$ws.Range("P4").Value = 2.006342455445734
$ws.Range("Q4").Value = '[1.2641844311742707, 2.7485004797171966]'
$ws.Range("R4").Value = [double]"2.064713407579077e-06"
$ws.Range("S4").Value = [double]"2.064713407579077e-06"
$ws.Range("T4").Value = 49.65087349391023
$ws.Range("U4").Value = '[35.765418302306614, 63.53632868551384]'
$ws.Range("V4").Value = [double]"5.118284684968444e-09"
$ws.Range("W4").Value = [double]"5.118284684968444e-09"
$ws.Range("X4").Value = 16.08448448448466
$ws.Range("Y4").Value = 13.2933533533535
$ws.Range("Z4").Value = 18.87561561561582

$ws.Range("F5").Value = 23.63000000000025
$ws.Range("H5").Value = 0.0008021299153820882
$ws.Range("I5").Value = 0.0008021299153820882
$ws.Range("L5").Value = 39.50255985554149
$ws.Range("M5").Value = '[14.636262894008638, 64.36885681707435]'
$ws.Range("N5").Value = 0.002522949339124114
$ws.Range("O5").Value = 0.002522949339124114
$ws.Range("P5").Value = 1.956026657190042
$ws.Range("Q5").Value = '[1.2641844311742725, 2.6478688832058106]'
$ws.Range("R5").Value = [double]"8.85456929777817e-07"
$ws.Range("S5").Value = [double]"8.85456929777817e-07"
$ws.Range("T5").Value = 70.75992198709156
$ws.Range("U5").Value = '[57.11666682869338, 84.40317714548974]'
$ws.Range("V5").Value = [double]"1.303401830909934e-13"
$ws.Range("W5").Value = [double]"1.303401830909934e-13"
$ws.Range("X5").Value = 16.27371371371389
$ws.Range("Y5").Value = 13.67181181181196
$ws.Range("Z5").Value = 18.87561561561581

$ws.Range("F6").Value = 23.63000000000025
$ws.Range("H6").Value = 0.009338072926212471
$ws.Range("I6").Value = 0.009338072926212471
$ws.Range("L6").Value = 38.97680370733997
$ws.Range("M6").Value = '[11.224914606551465, 66.72869280812847]'
$ws.Range("N6").Value = 0.006954648815170827
$ws.Range("O6").Value = 0.006954648815170827
$ws.Range("P6").Value = 1.66671081721981
$ws.Range("Q6").Value = '[0.6478159025420398, 2.6856057318975806]'
$ws.Range("R6").Value = 0.001926001307358982
$ws.Range("S6").Value = 0.001926001307358982
$ws.Range("T6").Value = 62.00758669462707
$ws.Range("U6").Value = '[44.92463554683647, 79.09053784241766]'
$ws.Range("V6").Value = [double]"3.533975112546273e-09"
$ws.Range("W6").Value = [double]"3.533975112546273e-09"
$ws.Range("X6").Value = 17.36178178178197
$ws.Range("Y6").Value = 13.52988988989004
$ws.Range("Z6").Value = 21.1936736736739

$ws.Range("F7").Value = 23.63000000000025
$ws.Range("H7").Value = 0.001904091200344293
$ws.Range("I7").Value = 0.001904091200344293
$ws.Range("L7").Value = 42.38546941194052
$ws.Range("M7").Value = '[15.13923678075247, 69.63170204312857]'
$ws.Range("N7").Value = 0.003039145272869526
$ws.Range("O7").Value = 0.003039145272869526
$ws.Range("P7").Value = 1.842816111114733
$ws.Range("Q7").Value = '[1.0629212381515023, 2.6227109840779645]'
$ws.Range("R7").Value = [double]"2.042131984247852e-05"
$ws.Range("S7").Value = [double]"2.042131984247852e-05"
$ws.Range("T7").Value = 65.08022410921416
$ws.Range("U7").Value = '[49.30785828838891, 80.85258993003941]'
$ws.Range("V7").Value = [double]"1.222504319997597e-10"
$ws.Range("W7").Value = [double]"1.222504319997597e-10"
$ws.Range("X7").Value = 16.69947947947966
$ws.Range("Y7").Value = 13.76642642642658
$ws.Range("Z7").Value = 19.63253253253274

$ws.Range("F8").Value = 23.63000000000025
$ws.Range("H8").Value = 0.0009370381424760321
$ws.Range("I8").Value = 0.0009370381424760321
$ws.Range("L8").Value = 35.65443238682776
$ws.Range("M8").Value = '[11.607211686352478, 59.701653087303036]'
$ws.Range("N8").Value = 0.004556801460608551
$ws.Range("O8").Value = 0.004556801460608551
$ws.Range("P8").Value = 2.182447749340657
$ws.Range("Q8").Value = '[1.515763422452732, 2.8491320762285817]'
$ws.Range("R8").Value = [double]"4.094951555622117e-08"
$ws.Range("S8").Value = [double]"4.094951555622117e-08"
$ws.Range("T8").Value = 56.45263197720858
$ws.Range("U8").Value = '[43.93886703060309, 68.96639692381407]'
$ws.Range("V8").Value = [double]"9.580558568700326e-12"
$ws.Range("W8").Value = [double]"9.580558568700326e-12"
$ws.Range("X8").Value = 15.42218218218235
$ws.Range("Y8").Value = 12.91489489489503
$ws.Range("Z8").Value = 17.92946946946967

$ws.Range("F9").Value = 23.63000000000025
$ws.Range("H9").Value = [double]"7.214132849320265e-06"
$ws.Range("I9").Value = [double]"7.214132849320265e-06"
$ws.Range("L9").Value = 58.05715937802492
$ws.Range("M9").Value = '[31.360638707300936, 84.7536800487489]'
$ws.Range("N9").Value = [double]"7.0082230603008e-05"
$ws.Range("O9").Value = [double]"7.0082230603008e-05"
$ws.Range("P9").Value = 2.03150035457358
$ws.Range("Q9").Value = '[1.553500271144502, 2.5095004380026573]'
$ws.Range("R9").Value = [double]"5.355205168200428e-11"
$ws.Range("S9").Value = [double]"5.355205168200428e-11"
$ws.Range("T9").Value = 68.13070013294325
$ws.Range("U9").Value = '[53.38989578389733, 82.87150448198916]'
$ws.Range("V9").Value = [double]"4.664491015660133e-12"
$ws.Range("W9").Value = [double]"4.664491015660133e-12"
$ws.Range("X9").Value = 15.98986986987004
$ws.Range("Y9").Value = 14.19219219219234
$ws.Range("Z9").Value = 17.78754754754774

$ws.Range("F10").Value = 23.63000000000025
$ws.Range("H10").Value = 0.003607841918666921
$ws.Range("I10").Value = 0.003607841918666921
$ws.Range("L10").Value = 41.81680708990871
$ws.Range("M10").Value = '[12.456980525391842, 71.17663365442557]'
$ws.Range("N10").Value = 0.006255487334577259
$ws.Range("O10").Value = 0.006255487334577259
$ws.Range("P10").Value = 1.86797401024258
$ws.Range("Q10").Value = '[1.0755001877154244, 2.6604478327697363]'
$ws.Range("R10").Value = [double]"2.121576289471072e-05"
$ws.Range("S10").Value = [double]"2.121576289471072e-05"
$ws.Range("T10").Value = 67.13583366543352
$ws.Range("U10").Value = '[50.59416844634006, 83.67749888452698]'
$ws.Range("V10").Value = [double]"1.923965431416264e-10"
$ws.Range("W10").Value = [double]"1.923965431416264e-10"
$ws.Range("X10").Value = 16.60486486486504
$ws.Range("Y10").Value = 13.62450450450465
$ws.Range("Z10").Value = 19.58522522522544

$ws.Range("F11").Value = 24.08000000000033
$ws.Range("H11").Value = 0.001650856289021396
$ws.Range("I11").Value = 0.001650856289021396
$ws.Range("L11").Value = 34.65303227317883
$ws.Range("M11").Value = '[10.197267223805632, 59.10879732255203]'
$ws.Range("N11").Value = 0.006505771481017497
$ws.Range("O11").Value = 0.006505771481017497
$ws.Range("P11").Value = 2.232763547596349
$ws.Range("Q11").Value = '[1.566079220708425, 2.899447874484273]'
$ws.Range("R11").Value = [double]"2.433559176395761e-08"
$ws.Range("S11").Value = [double]"2.433559176395761e-08"
$ws.Range("T11").Value = 50.53991019096253
$ws.Range("U11").Value = '[37.7149661098042, 63.36485427212086]'
$ws.Range("V11").Value = [double]"4.256837105032218e-10"
$ws.Range("W11").Value = [double]"4.256837105032218e-10"
$ws.Range("X11").Value = 15.52304304304325
$ws.Range("Y11").Value = 12.96800800800819
$ws.Range("Z11").Value = 18.07807807807832

$ws.Range("B12").Value = 1
$ws.Range("F12").Value = 24.08000000000033
$ws.Range("H12").Value = [double]"4.332876316703871e-05"
$ws.Range("I12").Value = [double]"4.332876316703871e-05"
$ws.Range("L12").Value = 48.28464312305041
$ws.Range("M12").Value = '[21.976595271173792, 74.59269097492702]'
$ws.Range("N12").Value = 0.0005910072927881593
$ws.Range("O12").Value = 0.0005910072927881593
$ws.Range("P12").Value = 2.396289891927349
$ws.Range("Q12").Value = '[1.8931319093704246, 2.8994478744842738]'
$ws.Range("R12").Value = [double]"1.885158695813516e-12"
$ws.Range("S12").Value = [double]"1.885158695813516e-12"
$ws.Range("T12").Value = 60.11912604889449
$ws.Range("U12").Value = '[46.46412149591835, 73.77413060187062]'
$ws.Range("V12").Value = [double]"1.951172556857728e-11"
$ws.Range("W12").Value = [double]"1.951172556857728e-11"
$ws.Range("X12").Value = 14.89633633633654
$ws.Range("Y12").Value = 12.96800800800818
$ws.Range("Z12").Value = 16.8246646646649

$ws.Range("F13").Value = 24.08000000000033
$ws.Range("H13").Value = 0.001103041852135833
$ws.Range("I13").Value = 0.001103041852135833
$ws.Range("L13").Value = 43.35688085424794
$ws.Range("M13").Value = '[14.025451859982326, 72.68830984851355]'
$ws.Range("N13").Value = 0.004670856531870493
$ws.Range("O13").Value = 0.004670856531870493
$ws.Range("P13").Value = 2.03150035457358
$ws.Range("Q13").Value = '[1.3270791789938867, 2.7359215301532727]'
$ws.Range("R13").Value = [double]"6.003055181835748e-07"
$ws.Range("S13").Value = [double]"6.003055181835748e-07"
$ws.Range("T13").Value = 62.09111257003961
$ws.Range("U13").Value = '[46.606507322561484, 77.57571781751774]'
$ws.Range("V13").Value = [double]"2.670292875706082e-10"
$ws.Range("W13").Value = [double]"2.670292875706082e-10"
$ws.Range("X13").Value = 16.29437437437459
$ws.Range("Y13").Value = 13.5947147147149
$ws.Range("Z13").Value = 18.99403403403429

$ws.Range("B14").Value = 0
$ws.Range("F14").Value = 24.08000000000033
$ws.Range("H14").Value = 0.036532143778503
$ws.Range("I14").Value = 0.036532143778503
$ws.Range("L14").Value = 27.72367581638446
$ws.Range("M14").Value = '[-0.9260649439240822, 56.37341657669301]'
$ws.Range("N14").Value = 0.05754567945071698
$ws.Range("O14").Value = 0.05754567945071698
$ws.Range("P14").Value = 2.383710942363427
$ws.Range("Q14").Value = '[0.421394810391424, 4.34602707433543]'
$ws.Range("R14").Value = 0.01839121727147575
$ws.Range("S14").Value = 0.01839121727147575
$ws.Range("T14").Value = 58.91722120609735
$ws.Range("U14").Value = '[44.138173012564536, 73.69626939963017]'
$ws.Range("V14").Value = [double]"3.124855929570458e-10"
$ws.Range("W14").Value = [double]"3.124855929570458e-10"
$ws.Range("X14").Value = 14.94454454454475
$ws.Range("Y14").Value = 7.424064064064165
$ws.Range("Z14").Value = 22.46502502502533
